# Weekly update: a new price record was added for "Feria Lagunitas de
# Puerto Montt - Cebollín" (Hortaliza). This pushes the existing
# historical rows down by one and introduces a fresh row with the most
# recent date, re-using the same Volumen/Precio/Unidad data as the row
# it displaces (rows shift down; the last historical row simply moves
# from row 431 to row 432 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 306 -- everything
# currently on rows 306:431 (including the last row, 431) shifts down
# to 307:432 automatically, carrying over values and formatting.
$ws.Rows.Item(306).Insert()

# Populate the newly inserted row 306 with the new weekly data point.
$ws.Cells.Item(306, 1).Value = 4
$ws.Cells.Item(306, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(306, 3).Value = "Los Lagos"
$ws.Cells.Item(306, 4).Value = 45027
$ws.Cells.Item(306, 5).Value = 10
$ws.Cells.Item(306, 6).Value = 100112037
$ws.Cells.Item(306, 7).Value = "Cebollín"
$ws.Cells.Item(306, 8).Value = "Sin especificar"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 180
$ws.Cells.Item(306, 11).Value = 6500
$ws.Cells.Item(306, 12).Value = 7000
$ws.Cells.Item(306, 13).Value = 6750
$ws.Cells.Item(306, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(306, 15).Value = "Región Metropolitana"
$ws.Cells.Item(306, 16).Value = 188
$ws.Cells.Item(306, 17).Value = 36
$ws.Cells.Item(306, 18).Value = "Hortaliza"
